# Moved RegistrationForm to forms and made necessary changes
# Adds a new row (12) of sample user data to the User sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 10.0
$ws.Range("B12").Value = "S7654321A"
$ws.Range("C12").Value = "Password1234"
